# Mejor descripción para la plantilla que se proporciona para el registro de materias
#
# The "Septimo" (no accent) entries used as the "ciclo" value for column B
# were replaced with a properly accented "Séptimo" so the template text
# reads correctly. This affects every row that referenced that cycle.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

foreach ($r in 4, 16, 26, 30, 37) {
    $cell = $ws.Cells.Item($r, 2)
    if ($cell.Text -eq "Septimo") {
        $cell.Value2 = "Séptimo"
    }
}

# Reflect where the author ended up navigating/selecting in the sheet.
$ws.Activate()
$aw = $excel.ActiveWindow
$aw.ScrollRow = 32
$aw.ScrollColumn = 1
$ws.Range("B37").Select()
